$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for the two new columns, mirroring the style of the existing
# header row (bold / bordered / centered), copied from column H's header.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Row-by-row data for the new I (I0) and J (IF) columns.
$data = @(
    @(2, 6, 6),
    @(3, 10, 10),
    @(4, 10, 10),
    @(5, 7, 7),
    @(6, 7, 7),
    @(7, 8, 8),
    @(8, 7, 7),
    @(9, 7, 7),
    @(10, 7, 7),
    @(11, 7, 8),
    @(12, 7, 7),
    @(13, 7, 7),
    @(14, 8, 8),
    @(15, 7, 7),
    @(16, 6, 7),
    @(17, 6, 7),
    @(18, 7, 7),
    @(19, 8, 8),
    @(20, 8, 8),
    @(21, 7, 7),
    @(22, 7, 7),
    @(23, 7, 7),
    @(24, 7, 7),
    @(25, 6, 7),
    @(26, 7, 7),
    @(27, 8, 8),
    @(28, 7, 7),
    @(29, 6, 7),
    @(30, 8, 8),
    @(31, 7, 7),
    @(32, 8, 8),
    @(33, 5, 7),
    @(34, 10, 10),
    @(35, 6, 6),
    @(36, 9, 9),
    @(37, 6, 6),
    @(38, 7, 7),
    @(39, 8, 8),
    @(40, 9, 9),
    @(41, 9, 9),
    @(42, 7, 7),
    @(43, 9, 9),
    @(44, 4, 5),
    @(45, 8, 8),
    @(46, 7, 7),
    @(47, 6, 7),
    @(48, 7, 7),
    @(49, 6, 6),
    @(50, 6, 6),
    @(51, 8, 8),
    @(52, 7, 7),
    @(53, 8, 8),
    @(54, 7, 7)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 9).Value = $entry[1]
    $ws.Cells.Item($r, 10).Value = $entry[2]
}
